# Duplicate the two existing data rows (2 and 3) into new rows 4 and 5,
# keeping the same values, row heights and styling (wrapText on column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 = copy of row 2
$ws.Range("A4").Value = $ws.Range("A2").Value2
$ws.Range("B4").Value = $ws.Range("B2").Value2
$ws.Range("C4").Value = $ws.Range("C2").Value2
$ws.Range("D4").Value = $ws.Range("D2").Value2
$ws.Rows.Item(4).RowHeight = $ws.Rows.Item(2).RowHeight
$ws.Range("D4").WrapText = $true

# Row 5 = copy of row 3
$ws.Range("A5").Value = $ws.Range("A3").Value2
$ws.Range("B5").Value = $ws.Range("B3").Value2
$ws.Range("C5").Value = $ws.Range("C3").Value2
$ws.Range("D5").Value = $ws.Range("D3").Value2
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(3).RowHeight
$ws.Range("D5").WrapText = $true
